$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Proximity sheet: append rows 43-48 (Living Room Main Door ENTER/EXIT events)
# ----------------------------------------------------------------------------
$proximity = $wb.Worksheets.Item("Proximity")

$pRowNums = @(43, 44, 45, 46, 47, 48)
$pColA = @("2026-02-01", "2026-02-01", "2026-02-01", "2026-02-01", "2026-02-01", "2026-02-01")
$pColB = @("14:23:00", "14:23:02", "14:23:25", "14:23:32", "14:23:44", "14:23:49")
$pColC = @("14:00", "14:00", "14:00", "14:00", "14:00", "14:00")
$pColD = @("Living Room Main Door", "Living Room Main Door", "Living Room Main Door", "Living Room Main Door", "Living Room Main Door", "Living Room Main Door")
$pColE = @("ENTER", "EXIT", "ENTER", "EXIT", "ENTER", "EXIT")
$pColF = @("User ENTERED Living Room Main Door", "User EXITED Living Room Main Door", "User ENTERED Living Room Main Door", "User EXITED Living Room Main Door", "User ENTERED Living Room Main Door", "User EXITED Living Room Main Door")

for ($i = 0; $i -lt $pRowNums.Count; $i++) {
    $rn = $pRowNums[$i]
    $proximity.Cells.Item($rn, 1).NumberFormat = "@"
    $proximity.Cells.Item($rn, 1).Value = $pColA[$i]
    $proximity.Cells.Item($rn, 2).NumberFormat = "@"
    $proximity.Cells.Item($rn, 2).Value = $pColB[$i]
    $proximity.Cells.Item($rn, 3).NumberFormat = "@"
    $proximity.Cells.Item($rn, 3).Value = $pColC[$i]
    $proximity.Cells.Item($rn, 4).NumberFormat = "@"
    $proximity.Cells.Item($rn, 4).Value = $pColD[$i]
    $proximity.Cells.Item($rn, 5).NumberFormat = "@"
    $proximity.Cells.Item($rn, 5).Value = $pColE[$i]
    $proximity.Cells.Item($rn, 6).NumberFormat = "@"
    $proximity.Cells.Item($rn, 6).Value = $pColF[$i]
}

# ----------------------------------------------------------------------------
# mmWave sheet: append rows 2-6 (Living Room presence events)
# ----------------------------------------------------------------------------
$mmwave = $wb.Worksheets.Item("mmWave")

$mRowNums = @(2, 3, 4, 5, 6)
$mColA = @("2026-02-01", "2026-02-01", "2026-02-01", "2026-02-01", "2026-02-01")
$mColB = @("14:23:12", "14:23:23", "14:23:35", "14:23:43", "14:23:54")
$mColC = @("14:00", "14:00", "14:00", "14:00", "14:00")
$mColD = @("Living Room", "Living Room", "Living Room", "Living Room", "Living Room")
$mColE = @("NO_MOTION_DETECTED", "PRESENCE_DETECTED", "PRESENCE_DETECTED", "PRESENCE_DETECTED", "PRESENCE_DETECTED")
$mColF = @("Inactive", "Active", "Active", "Active", "Active")

for ($i = 0; $i -lt $mRowNums.Count; $i++) {
    $rn = $mRowNums[$i]
    $mmwave.Cells.Item($rn, 1).NumberFormat = "@"
    $mmwave.Cells.Item($rn, 1).Value = $mColA[$i]
    $mmwave.Cells.Item($rn, 2).NumberFormat = "@"
    $mmwave.Cells.Item($rn, 2).Value = $mColB[$i]
    $mmwave.Cells.Item($rn, 3).NumberFormat = "@"
    $mmwave.Cells.Item($rn, 3).Value = $mColC[$i]
    $mmwave.Cells.Item($rn, 4).NumberFormat = "@"
    $mmwave.Cells.Item($rn, 4).Value = $mColD[$i]
    $mmwave.Cells.Item($rn, 5).NumberFormat = "@"
    $mmwave.Cells.Item($rn, 5).Value = $mColE[$i]
    $mmwave.Cells.Item($rn, 6).NumberFormat = "@"
    $mmwave.Cells.Item($rn, 6).Value = $mColF[$i]
}

# ----------------------------------------------------------------------------
# Camera sheet: append rows 26-30 (Living Room Main Door image captures)
# ----------------------------------------------------------------------------
$camera = $wb.Worksheets.Item("Camera")

$cRowNums = @(26, 27, 28, 29, 30)
$cColA = @("2026-02-01", "2026-02-01", "2026-02-01", "2026-02-01", "2026-02-01")
$cColB = @("14:23:02", "14:23:31", "14:23:34", "14:23:48", "14:23:53")
$cColC = @("14:00", "14:00", "14:00", "14:00", "14:00")
$cColD = @("Living Room Main Door", "Living Room Main Door", "Living Room Main Door", "Living Room Main Door", "Living Room Main Door")
$cColE = @("Image Captured", "Image Captured", "Image Captured", "Image Captured", "Image Captured")
$cColF = @("Active", "Active", "Active", "Active", "Active")

for ($i = 0; $i -lt $cRowNums.Count; $i++) {
    $rn = $cRowNums[$i]
    $camera.Cells.Item($rn, 1).NumberFormat = "@"
    $camera.Cells.Item($rn, 1).Value = $cColA[$i]
    $camera.Cells.Item($rn, 2).NumberFormat = "@"
    $camera.Cells.Item($rn, 2).Value = $cColB[$i]
    $camera.Cells.Item($rn, 3).NumberFormat = "@"
    $camera.Cells.Item($rn, 3).Value = $cColC[$i]
    $camera.Cells.Item($rn, 4).NumberFormat = "@"
    $camera.Cells.Item($rn, 4).Value = $cColD[$i]
    $camera.Cells.Item($rn, 5).NumberFormat = "@"
    $camera.Cells.Item($rn, 5).Value = $cColE[$i]
    $camera.Cells.Item($rn, 6).NumberFormat = "@"
    $camera.Cells.Item($rn, 6).Value = $cColF[$i]
}
